{"js": "// Optimizaci\u00f3n de lista de asistencia\n// - Title: add \"zona 7\" suffix\n// - Fecha: 2025-02-15 -> 2025-02-16\n// - Header row: disambiguate the two \"Firma\" columns -> \"Firma Entrada\" / \"Firma Salida\"\n// - Row 2 (AGABO MARTINEZ MIGUEL ANGEL): Hora de Entrada/Salida updated\n// - Row 3: name replaced (ALBARRAN JIMENEZ -> ALBARRAN GARCIA) and both times updated\n\nconst body = context.document.body;\n\n// --- Title paragraph ---\nconst titleResults = body.search(\"Lista de asistencia: concurso de escoltas\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\ntitleResults.items[0].insertText(\"Lista de asistencia: concurso de escoltas zona 7\", \"Replace\");\n\n// --- Date paragraph ---\nconst dateResults = body.search(\"Fecha: 2025-02-15\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\ndateResults.items[0].insertText(\"Fecha: 2025-02-16\", \"Replace\");\n\nawait context.sync();\n\n// --- Table edits (first table in the document, 0-indexed rows/cols) ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Header row: two \"Firma\" cells -> \"Firma Entrada\" / \"Firma Salida\"\ntable.getCell(0, 3).value = \"Firma Entrada\";\ntable.getCell(0, 5).value = \"Firma Salida\";\n\n// Row index 1 (AGABO MARTINEZ MIGUEL ANGEL): entrada/salida hours\ntable.getCell(1, 2).value = \"04:45:00\";\ntable.getCell(1, 4).value = \"15:00:00\";\n\n// Row index 2: name + entrada/salida hours\ntable.getCell(2, 1).value = \"ALBARRAN  GARCIA  CARLOS ARTURO\";\ntable.getCell(2, 2).value = \"04:46:00\";\ntable.getCell(2, 4).value = \"15:00:00\";\n\nawait context.sync();\n", "ps1": "# Optimizaci\u00f3n de lista de asistencia\n# - Title: add \"zona 7\" suffix\n# - Fecha: 2025-02-15 -> 2025-02-16\n# - Header row: disambiguate the two \"Firma\" columns -> \"Firma Entrada\" / \"Firma Salida\"\n# - Row 2 (AGABO MARTINEZ MIGUEL ANGEL): Hora de Entrada/Salida updated\n# - Row 3: name replaced (ALBARRAN JIMENEZ -> ALBARRAN GARCIA) and both times updated\n\n$d = $word.ActiveDocument\n\n# --- Title / date paragraphs ---\n$find = $d.Content.Find\n$find.Execute(\"Lista de asistencia: concurso de escoltas\", $false, $false, $false, $false, $false, $true, 1, $false, \"Lista de asistencia: concurso de escoltas zona 7\", 1)\n\n$find = $d.Content.Find\n$find.Execute(\"Fecha: 2025-02-15\", $false, $false, $false, $false, $false, $true, 1, $false, \"Fecha: 2025-02-16\", 1)\n\n# --- Table edits (first table in the document) ---\n$t = $d.Tables(1)\n\n# Header row: two \"Firma\" cells -> \"Firma Entrada\" / \"Firma Salida\"\n$t.Cell(1, 4).Range.Text = \"Firma Entrada\"\n$t.Cell(1, 6).Range.Text = \"Firma Salida\"\n\n# Row 2 (AGABO MARTINEZ MIGUEL ANGEL): entrada/salida hours\n$t.Cell(2, 3).Range.Text = \"04:45:00\"\n$t.Cell(2, 5).Range.Text = \"15:00:00\"\n\n# Row 3: name + entrada/salida hours\n$t.Cell(3, 2).Range.Text = \"ALBARRAN  GARCIA  CARLOS ARTURO\"\n$t.Cell(3, 3).Range.Text = \"04:46:00\"\n$t.Cell(3, 5).Range.Text = \"15:00:00\"\n"}
